$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old hyperlink currently on B3 (Lolo's email cell) before
# shifting columns, so we don't end up with a stale relationship.
$ws.Range("B3").Hyperlinks.Delete()

# Insert a new column A (shifts existing A,B,C -> B,C,D) to make room
# for a row-number / index column.
$ws.Columns("A").Insert()

# Fill the new index column (0-based row counter for the 3 data rows).
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2

# Update the two email addresses that changed.
$ws.Range("C3").Value = "hans.che@donebyngle.com"
$ws.Range("C4").Value = "sxdp3fpkzql@ezztt.com"

# Re-create the hyperlink on the (now shifted) Lolo email cell.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:hans.che@donebyngle.com", "", "", "hans.che@donebyngle.com")

# Style the header row (bold, thin box border, centered/top aligned).
$headerRange = $ws.Range("B1:D1")
$headerRange.Font.Bold = $true
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 11
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Apply the same style to the new index column's data cells.
$indexRange = $ws.Range("A2:A4")
$indexRange.Font.Bold = $true
$indexRange.Font.Name = "Calibri"
$indexRange.Font.Size = 11
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# Match the page setup / selection tweaks from the edit.
$ws.PageSetup.Orientation = 1
$ws.Range("C3").Select()
